$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.06091455863883025
$ws.Range("H2").Value = -8.147546183989046
$ws.Range("I2").Value = -5.123253149119593
$ws.Range("G3").Value = 0.1307271089491298
$ws.Range("H3").Value = 10.53942306471695
$ws.Range("G4").Value = -0.4145087063810201
$ws.Range("H4").Value = -50.73298695194668
$ws.Range("G5").Value = -0.4290888256076349
$ws.Range("H5").Value = -7.537052995745741
$ws.Range("G6").Value = 0.1836515653084748
$ws.Range("H6").Value = -6.8463973129281
$ws.Range("G7").Value = 0.3226415608039043
$ws.Range("H7").Value = 55.57817252010035
$ws.Range("G8").Value = 0.1015084768147391
$ws.Range("H8").Value = -0.3803294438824816
$ws.Range("G9").Value = 0.1370387726508655
$ws.Range("H9").Value = 8.34812326071417
$ws.Range("G10").Value = 0.04704933322488328
$ws.Range("H10").Value = -23.4195729607367
$ws.Range("G11").Value = 0.03773178224432856
$ws.Range("H11").Value = -24.43086708821197
$ws.Range("G12").Value = 0.1219543408696468
$ws.Range("H12").Value = 31.74680235837442
$ws.Range("G13").Value = 0.08924474581767286
$ws.Range("H13").Value = 17.10722664915843
$ws.Range("G14").Value = 0.231658027129248
$ws.Range("H14").Value = 2.52095176261215
$ws.Range("G15").Value = 0.2457655062503053
$ws.Range("H15").Value = -0.2423496092122862
$ws.Range("G16").Value = 0.119046778688261
$ws.Range("H16").Value = 4.661239920661622
$ws.Range("G17").Value = 0.1156884937184325
$ws.Range("H17").Value = -22.56894165611904
$ws.Range("G18").Value = -0.01068220093323061
$ws.Range("H18").Value = -19.3282126703506
$ws.Range("G19").Value = 0.0425761115824068
$ws.Range("H19").Value = 75.77607997317403
$ws.Range("G20").Value = 0.1244012655429656
$ws.Range("H20").Value = 46.24920040055695
$ws.Range("G21").Value = 0.1338669105052991
$ws.Range("H21").Value = 104.5207010044334
$ws.Range("G22").Value = 0.1835804114812105
$ws.Range("H22").Value = -4.166421649889425
$ws.Range("G23").Value = 0.2069301937187523
$ws.Range("H23").Value = -4.068999163984525
$ws.Range("G24").Value = -0.01935087340486509
$ws.Range("H24").Value = -408.4883566047319
$ws.Range("G25").Value = 0.003654445474025173
$ws.Range("H25").Value = 115.7130971288975
$ws.Range("G26").Value = 0.185242108368188
$ws.Range("H26").Value = -9.579033223270891
$ws.Range("G27").Value = 0.2048634857006695
$ws.Range("H27").Value = 6.210422682659293
$ws.Range("G28").Value = 0.04606343700052995
$ws.Range("H28").Value = -31.15990370764158
$ws.Range("G29").Value = 0.09301504621873845
$ws.Range("H29").Value = -1.325557678766941

Write-Output "Updated cells with new compared-return values."
